# Adds ODATA dev and data generation tasks to the Tasks sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Task numbers (column A)
$ws.Cells.Item(20, 1).Value = 17
$ws.Cells.Item(21, 1).Value = 18
$ws.Cells.Item(22, 1).Value = 19

# Descriptions (column B) for the two ODATA rows, entered together first.
$ws.Cells.Item(20, 2).Value = "ODATA. Add navigation between solution category and solutions"
$ws.Cells.Item(21, 2).Value = "ODATA. Implement substringof funtion for solution description field"

# Responsible (column C)
$ws.Cells.Item(20, 3).Value = "Shamil"
$ws.Cells.Item(21, 3).Value = "Shamil"

# Done criteria (column D)
$ws.Cells.Item(20, 4).Value = "Task 2 is finished"
$ws.Cells.Item(21, 4).Value = "Task 2 is finished"

# Prerequisites (column E)
$ws.Cells.Item(20, 5).Value = "test"
$ws.Cells.Item(21, 5).Value = "test"

# Status (column F)
$ws.Cells.Item(20, 6).Value = "Done"
$ws.Cells.Item(21, 6).Value = "Done"

# Time (column G)
$ws.Cells.Item(20, 7).Value = "0,5 day"
$ws.Cells.Item(21, 7).Value = "0,5 day"

# Row 22: Generate test data. Upload data from Historical data xls
$ws.Cells.Item(22, 2).Value = "Generate test data. Upload data from Historical data xls"
$ws.Cells.Item(22, 3).Value = "Shamil"
$ws.Cells.Item(22, 4).Value = "Task 12 is finished"
$ws.Cells.Item(22, 5).Value = "test"
$ws.Cells.Item(22, 6).Value = "In process"
$ws.Cells.Item(22, 7).Value = "1 hour"

# Copy formatting (wrap text, fills) from existing rows so the new rows match the
# established table style.
$ws.Range("A20:G21").WrapText = $true
$ws.Range("F20:F21").Interior.Color = $ws.Range("F2").Interior.Color

$ws.Range("A22:E22").WrapText = $true
$ws.Range("F22").Interior.Color = $ws.Range("F10").Interior.Color

# Update the visible selection/view to match the end of the edit session.
$ws.Range("D22").Select()
$excel.ActiveWindow.ScrollRow = 13
